$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.135837078094482
$ws.Range("B1").Value = 2.832927465438843
$ws.Range("C1").Value = 3.880904674530029
$ws.Range("D1").Value = 3.709670305252075
$ws.Range("E1").Value = 1.226359367370605
